{"js": "// Append new sentences about making natural weed sprays and a note\n// about strawberries to the end of the notes paragraph (the one that\n// currently ends with \"...alternatives for weeds.\").\nconst body = context.document.body;\n\nconst addition =\n  \" You can even look up online how to make these sprays. Strawberries  \";\n\nconst results = body.search(\"alternatives for weeds.\", { matchCase: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Anchor on the located text and grow to its containing paragraph so the\n  // new text lands at the very end of that paragraph (a fresh run),\n  // regardless of where the paragraph happens to sit in the document.\n  const para = results.items[0].paragraphs.getFirst();\n  para.insertText(addition, Word.InsertLocation.end);\n} else {\n  // Fallback: append to the end of the last paragraph of the document.\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n  const last = paragraphs.items[paragraphs.items.length - 1];\n  last.insertText(addition, Word.InsertLocation.end);\n}\n\nawait context.sync();\n", "ps1": "# Append new sentences about making natural weed sprays and a note\n# about strawberries to the end of the notes paragraph (the one that\n# currently ends with \"...alternatives for weeds.\").\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\"alternatives for weeds.\")\n\nif ($found) {\n    $rng.Collapse(0)  # wdCollapseEnd\n    $rng.InsertAfter(\" You can even look up online how to make these sprays. Strawberries  \")\n} else {\n    # Fallback: append to the end of the last paragraph of the document.\n    $p = $d.Paragraphs($d.Paragraphs.Count)\n    $r = $p.Range\n    $r.End = $r.End - 1\n    $r.InsertAfter(\" You can even look up online how to make these sprays. Strawberries  \")\n}\n"}
